$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2084942084942085
$ws.Range("C2").Value = 0.5366795366795367
$ws.Range("J2").Value = 0.007722007722007722
$ws.Range("P2").Value = 0.1505791505791506
$ws.Range("S2").Value = 0.09652509652509653
$ws.Range("B3").Value = 0.00684931506849315
$ws.Range("C3").Value = 0.02054794520547945
$ws.Range("J3").Value = 0.0136986301369863
$ws.Range("P3").Value = 0.6438356164383562
$ws.Range("S3").Value = 0.3150684931506849
$ws.Range("J4").Value = 0.02702702702702703
$ws.Range("P4").Value = 0.6486486486486487
$ws.Range("S4").Value = 0.3243243243243243
$ws.Range("B6").Value = 0.02403846153846154
$ws.Range("D6").Value = 0.02884615384615385
$ws.Range("F6").Value = 0.0576923076923077
$ws.Range("J6").Value = 0.25
$ws.Range("O6").Value = 0.02884615384615385
$ws.Range("Q6").Value = 0.1346153846153846
$ws.Range("R6").Value = 0.08173076923076923
$ws.Range("S6").Value = 0.3942307692307692
$ws.Range("B7").Value = 0.1223404255319149
$ws.Range("D7").Value = 0.01063829787234043
$ws.Range("F7").Value = 0.04787234042553191
$ws.Range("J7").Value = 0.1223404255319149
$ws.Range("O7").Value = 0.01595744680851064
$ws.Range("Q7").Value = 0.1542553191489362
$ws.Range("R7").Value = 0.07446808510638298
$ws.Range("S7").Value = 0.4521276595744681
$ws.Range("B8").Value = 0.091324200913242
$ws.Range("D8").Value = 0.01141552511415525
$ws.Range("E8").Value = 0.00228310502283105
$ws.Range("F8").Value = 0.045662100456621
$ws.Range("J8").Value = 0.1164383561643836
$ws.Range("O8").Value = 0.0091324200913242
$ws.Range("Q8").Value = 0.1712328767123288
$ws.Range("R8").Value = 0.09817351598173515
$ws.Range("S8").Value = 0.454337899543379
$ws.Range("B9").Value = 0.0861244019138756
$ws.Range("D9").Value = 0.01913875598086124
$ws.Range("F9").Value = 0.05741626794258373
$ws.Range("J9").Value = 0.1196172248803828
$ws.Range("O9").Value = 0.02392344497607655
$ws.Range("Q9").Value = 0.138755980861244
$ws.Range("R9").Value = 0.07655502392344497
$ws.Range("S9").Value = 0.4784688995215311
$ws.Range("B10").Value = 0.09127625201938611
$ws.Range("D10").Value = 0.01857835218093699
$ws.Range("E10").Value = 0.002423263327948304
$ws.Range("F10").Value = 0.07512116316639741
$ws.Range("J10").Value = 0.1195476575121163
$ws.Range("O10").Value = 0.01857835218093699
$ws.Range("Q10").Value = 0.197092084006462
$ws.Range("R10").Value = 0.09854604200323101
$ws.Range("S10").Value = 0.3788368336025848
$ws.Range("G11").Value = 0.1318327974276527
$ws.Range("J11").Value = 0.09967845659163987
$ws.Range("K11").Value = 0.1929260450160772
$ws.Range("L11").Value = 0.5530546623794212
$ws.Range("S11").Value = 0.022508038585209
$ws.Range("G12").Value = 0.7247191011235955
$ws.Range("J12").Value = 0.2134831460674157
$ws.Range("K12").Value = 0.005617977528089887
$ws.Range("L12").Value = 0.02247191011235955
$ws.Range("S12").Value = 0.03370786516853932
$ws.Range("G13").Value = 0.5476190476190477
$ws.Range("J13").Value = 0.3809523809523809
$ws.Range("S13").Value = 0.07142857142857142
$ws.Range("F15").Value = 0.01047120418848168
$ws.Range("H15").Value = 0.1308900523560209
$ws.Range("I15").Value = 0.06806282722513089
$ws.Range("J15").Value = 0.3298429319371728
$ws.Range("K15").Value = 0.05235602094240838
$ws.Range("M15").Value = 0.02094240837696335
$ws.Range("O15").Value = 0.07329842931937172
$ws.Range("S15").Value = 0.3141361256544503
$ws.Range("F16").Value = 0.01935483870967742
$ws.Range("H16").Value = 0.1419354838709677
$ws.Range("I16").Value = 0.08387096774193549
$ws.Range("J16").Value = 0.4193548387096774
$ws.Range("K16").Value = 0.1419354838709677
$ws.Range("M16").Value = 0.01290322580645161
$ws.Range("O16").Value = 0.03870967741935484
$ws.Range("S16").Value = 0.1419354838709677
$ws.Range("F17").Value = 0.01240694789081886
$ws.Range("H17").Value = 0.1315136476426799
$ws.Range("I17").Value = 0.08933002481389578
$ws.Range("J17").Value = 0.4491315136476427
$ws.Range("K17").Value = 0.09925558312655088
$ws.Range("M17").Value = 0.01488833746898263
$ws.Range("N17").Value = 0.002481389578163772
$ws.Range("O17").Value = 0.05955334987593052
$ws.Range("S17").Value = 0.141439205955335
$ws.Range("F18").Value = 0.01421800947867299
$ws.Range("H18").Value = 0.1848341232227488
$ws.Range("I18").Value = 0.08530805687203792
$ws.Range("J18").Value = 0.3554502369668247
$ws.Range("K18").Value = 0.1279620853080569
$ws.Range("M18").Value = 0.009478672985781991
$ws.Range("O18").Value = 0.06635071090047394
$ws.Range("S18").Value = 0.1563981042654028
$ws.Range("F19").Value = 0.01411589895988113
$ws.Range("H19").Value = 0.2258543833580981
$ws.Range("I19").Value = 0.09806835066864784
$ws.Range("J19").Value = 0.3580980683506687
$ws.Range("K19").Value = 0.1099554234769688
$ws.Range("M19").Value = 0.02228826151560178
$ws.Range("N19").Value = 0.001485884101040119
$ws.Range("O19").Value = 0.05349182763744428
$ws.Range("S19").Value = 0.1166419019316493
